# Applies the "Finished first two parts, need to add class diagrams and
# finish last part" edit to the technical-debt reflection document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Rewrite paragraph a) — drop the highlighted "IDK" placeholder and
#    replace with the fleshed-out technical-debt discussion.
# ---------------------------------------------------------------------
$paraA = $d.Paragraphs(3)
$rA = $paraA.Range
$rA.End = $paraA.Range.End - 1
$rA.Text = "a) Technical debt in our project can be defined in an iterative process as the debt associated with cutting corners to produce an iterative step. The technical debt that we have is when we create our squares, the dots are drawn over each other. This impacts the generator class and code that we have created. This is ok right now as we only need to draw the dots, but with later steps we must begin to average out the dot" + [char]0x2019 + "s color to create the color for the segments. Having multiple dots means that the color will be completely wrong and our code will need some major fixing. Having the wrong color means that our visualizer code needs to be fixed to get the right colour for our dots."

# ---------------------------------------------------------------------
# 2) Rewrite paragraph b) — expand with the extra clauses describing the
#    hash-table -> list data-structure change.
# ---------------------------------------------------------------------
$paraB = $d.Paragraphs(4)
$rB = $paraB.Range
$rB.End = $paraB.Range.End - 1
$rB.Text = "b) If we were to stay at the immutable data structure level, it would stop us from implementing a lot of functions in the future that would be much easier to do so if we were to change it up. For example, we found that the dots stored as hash tables were not too effective when drawing the segments because it all became very random and there was no order to it. We changed the data structure to be a list as we found having more order in our segments would help us in finishing our first step with business logic. If we were to have an immutable data structure, it would be incredibly difficult to try and work around it instead of just chancing our data structure. "

# ---------------------------------------------------------------------
# 3) Append PART 2 heading material, then the three list items, then the
#    PART 3 block — all as plain paragraphs first so nothing inherits
#    list formatting prematurely; numbering/style gets applied to the
#    three list paragraphs afterwards as one pass.
# ---------------------------------------------------------------------
$cur = $d.Paragraphs($d.Paragraphs.Count)
$cur.Range.Collapse(0)
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($d.Paragraphs.Count)
$cur.Range.Text = "PART 2"

$cur.Range.Collapse(0)
$cur.Range.InsertParagraphAfter()
$cur = $d.Paragraphs($d.Paragraphs.Count)
$cur.Range.Text = "WE NEED TO DRAW A CLASS DIAGRAM"

$cur.Range.Collapse(0)
$cur.Range.InsertParagraphAfter()
$list1 = $d.Paragraphs($d.Paragraphs.Count)
$list1.Range.Text = "To fix some of our flaws with the starter code, we first changed our data structure that would store each dot. It was a hash table. With has hash tables there is no order to them. It only has two corresponding values and the tables keys can be in any order. This made drawing the segments impossible and forced us to change everything to an array list. Having an array list allowed for a much better data structure to store our dots and segments. Likewise, we fixed up the problem with having multiple dots being drawn on each other by completely altering the code that was provided and adding if statements for some error that was occurring when drawing the segments."

$list1.Range.Collapse(0)
$list1.Range.InsertParagraphAfter()
$list2 = $d.Paragraphs($d.Paragraphs.Count)
$list2.Range.Text = "To support requests from the user we made everything a variable that can be easily changed in the command line argument. Likewise, utilized object orient programming to make it sure that if a class was not needed in our function, we were able to just not utilize the class and our code would not be impacted because of the nature of our composition in our code. Since only the classes that can be affected by invariants use composition, this means that our code can be manipulated and changed accordingly."

$list2.Range.Collapse(0)
$list2.Range.InsertParagraphAfter()
$list3 = $d.Paragraphs($d.Paragraphs.Count)
$list3.Range.Text = "The testing activity that we used helped us a lot. Mainly we used a lot of our testing to check where and how we could have problems arise in our code. However, we didn" + [char]0x2019 + "t utilize it enough when trying out our code, but this means that we can learn from our error and hopefully utilize it more in the final step of our code. Mainly, we used it to narrow down where the error in our code was and used it to narrow down the area that we needed to work on. Overall, we came to learn the importance of testing and we are now changing our mistakes and testing cases that can occur and other problems that may occur."

# ---------------------------------------------------------------------
# 4) Append a blank spacer + PART 3 heading material, all indented
#    (0.25in / 360 twips) rather than styled as list items.
# ---------------------------------------------------------------------
$list3.Range.Collapse(0)
$list3.Range.InsertParagraphAfter()
$blank1 = $d.Paragraphs($d.Paragraphs.Count)
$blank1.LeftIndent = 18

$blank1.Range.Collapse(0)
$blank1.Range.InsertParagraphAfter()
$part3 = $d.Paragraphs($d.Paragraphs.Count)
$part3.Range.Text = "PART 3"
$part3.LeftIndent = 18

$part3.Range.Collapse(0)
$part3.Range.InsertParagraphAfter()
$drawClass = $d.Paragraphs($d.Paragraphs.Count)
$drawClass.Range.Text = "DRAW A CLASS DIAGRAM"
$drawClass.LeftIndent = 18

$drawClass.Range.Collapse(0)
$drawClass.Range.InsertParagraphAfter()
$blank2 = $d.Paragraphs($d.Paragraphs.Count)
$blank2.LeftIndent = 18

# ---------------------------------------------------------------------
# 5) Turn the three PART 2 paragraphs into a single lowercase-letter
#    numbered list ("a)", "b)", "c)") using the List Paragraph style.
# ---------------------------------------------------------------------
$listRange = $d.Range($list1.Range.Start, $list3.Range.End)
$listRange.Style = "List Paragraph"
$listRange.ListFormat.ApplyNumberDefault()

$lvl = $list1.Range.ListFormat.ListTemplate.ListLevels(1)
$lvl.NumberStyle = 4
$lvl.NumberFormat = "%1)"

$st = $d.Styles("List Paragraph")
$st.Priority = 34
$st.NoSpaceBetweenParagraphsOfSameStyle = $true
$st.ParagraphFormat.LeftIndent = 36

Write-Host "Final paragraph count:" $d.Paragraphs.Count
